# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-01-11 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-01-12 Friday", 2)

# Update the 25 answer cells in the table (5 data rows x 5 columns,
# each data row followed by 3 blank rows). Cell.Range.Text assignment
# preserves the existing run formatting (font/size) of the cell.
$t = $d.Tables.Item(1)

$answers = @(
    "93÷9=10, 3", "90÷8=11, 2", "43÷8=5, 3", "80÷2=40, 0", "30÷3=10, 0",
    "48÷2=24, 0", "80÷7=11, 3", "10÷3=3, 1", "38÷2=19, 0", "21÷7=3, 0",
    "53÷3=17, 2", "87÷5=17, 2", "53÷2=26, 1", "85÷3=28, 1", "85÷9=9, 4",
    "62÷9=6, 8", "91÷6=15, 1", "40÷2=20, 0", "41÷6=6, 5", "12÷7=1, 5",
    "20÷7=2, 6", "45÷7=6, 3", "23÷2=11, 1", "43÷8=5, 3", "25÷9=2, 7"
)

$dataRows = @(1, 5, 9, 13, 17)

$idx = 0
foreach ($r in $dataRows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $answers[$idx]
        $idx = $idx + 1
    }
}
